# Sprint 3 backlog update (commit: "atualizacao do backlog da sprint 3")
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Backlog sprint 3": new rows for the HTML/CSS filters task and the
# "BD ou CVS" decision task, a new "legal name" task, row-height tweaks, and
# re-priced priorities that ripple from the new rows.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Backlog sprint 3")

# Capture donor cells for the priority styling (Alta / Média / Baixa) BEFORE
# any values change, so PasteSpecial always copies the right look.
$altaCell = $ws3.Range("E3")
$mediaCell = $ws3.Range("E6")
$baixaCell = $ws3.Range("E8")

# Write final cell contents for rows 3-12 (rows 5-12 replace/shift the old
# rows 5-10 to make room for the 3 new tasks).
$ws3.Range("A3").Value = "Estudo de Python3"
$ws3.Range("B3").Value = "Estudo de como será feita a integração/consulta do banco de dados com o site através do Python3"

$ws3.Range("A4").Value = "Integração com o site"
$ws3.Range("B4").Value = "Fazer o necessário para integrar/consultar o banco de dados, pelo site, através dos Python3"

$ws3.Range("A5").Value = "HTML filtros"
$ws3.Range("B5").Value = "Implementar a página de filtros no site"
$ws3.Range("C5").Value = "x"
$ws3.Range("D5").Value = "x"
$ws3.Range("E5").Value = "Alta"

$ws3.Range("A6").Value = "CSS filtros"
$ws3.Range("B6").Value = "Desenvolver o CSS para estilizar a pág. Filtros"
$ws3.Range("C6").Value = "x"
$ws3.Range("D6").Value = "x"
$ws3.Range("E6").Value = "Alta"

$ws3.Range("A7").Value = "BD ou CVS"
$ws3.Range("B7").Value = "perguntar se vale a pena trabalhar copm BD ou CSV"
$ws3.Range("C7").Value = "x"
$ws3.Range("D7").Value = "x"
$ws3.Range("E7").Value = "Média"

$ws3.Range("A8").Value = "Pesquisa medicamentos"
$ws3.Range("B8").Value = "Fazer as pesquisas sobre os medicamentos apontados na reunião com jornalista"
$ws3.Range("C8").Value = "x"
$ws3.Range("D8").Value = "x"
$ws3.Range("E8").Value = "Média"

$ws3.Range("A9").Value = "Pesquisa sintomas da covid longa e tratamentos"
$ws3.Range("B9").Value = "Procurar os sintomas pós-covid e como são tratados(cirurgias, uso de remedios, etc)"
$ws3.Range("C9").Value = "x"
$ws3.Range("D9").Value = "x"
$ws3.Range("E9").Value = "Média"

$ws3.Range("A10").Value = "Pesquisa investimentos"
$ws3.Range("B10").Value = "Procurar mais fontes sobre e investimentos, na área da saúde, de cada estado escolhido para o projeto(ou a falta deles)"
$ws3.Range("C10").Value = "x"
$ws3.Range("D10").Value = "x"
$ws3.Range("E10").Value = "Baixa"

$ws3.Range("A11").Value = "Responsividade "
$ws3.Range("B11").Value = "Atualizar o CSS para que se adapte para qualquer tamanho de tela sem quebrar"
$ws3.Range("C11").Value = "x"
$ws3.Range("D11").Value = "x"
$ws3.Range("E11").Value = "Baixa"

$ws3.Range("A12").Value = "Nome legal para o site"
$ws3.Range("B12").Value = "Criar um nome legal pro site e refazer a logo"
$ws3.Range("C12").Value = "x"
$ws3.Range("D12").Value = "x"
$ws3.Range("E12").Value = "Baixa"

# Re-apply the priority formatting (fill/font/border) that belongs to each
# bucket, since rows shifted and picked up new text.
$altaCell.Copy()
$ws3.Range("E4").PasteSpecial(-4122)
$ws3.Range("E5").PasteSpecial(-4122)
$ws3.Range("E6").PasteSpecial(-4122)

$mediaCell.Copy()
$ws3.Range("E7").PasteSpecial(-4122)
$ws3.Range("E8").PasteSpecial(-4122)
$ws3.Range("E9").PasteSpecial(-4122)

$baixaCell.Copy()
$ws3.Range("E10").PasteSpecial(-4122)
$ws3.Range("E11").PasteSpecial(-4122)
$ws3.Range("E12").PasteSpecial(-4122)

$ws3.Application.CutCopyMode = $false

# Row heights (auto-fit-style tweaks baked in by the author while editing).
$ws3.Rows.Item(3).RowHeight = 69
$ws3.Rows.Item(4).RowHeight = 68.25
$ws3.Rows.Item(5).RowHeight = 61.5
$ws3.Rows.Item(6).RowHeight = 61.5
$ws3.Rows.Item(7).RowHeight = 61.5
$ws3.Rows.Item(8).RowHeight = 61.5
$ws3.Rows.Item(9).RowHeight = 61.5
$ws3.Rows.Item(10).RowHeight = 61.5
$ws3.Rows.Item(11).RowHeight = 63
$ws3.Rows.Item(12).RowHeight = 63

# Scroll position / selection the author left the sheet in.
$ws3.Application.Goto($ws3.Range("A3"))
$ws3.Range("F12").Select()

# ---------------------------------------------------------------------------
# Sheet "Backlog sprint 4": the placeholder "TESTES" task is replaced by the
# real Python3-integration-study task (mirrors sprint 3's row 3), plus a new
# "Pesquisa sobre outras cidades" task and a "Tempo para conclusão" column.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Backlog sprint 4")

$ws4.Range("D1").Value = "Tempo para conclusão"
$ws4.Range("E1").Value = "Prioridade"

$ws4.Range("A2").Value = "Estudo de Python3"
$ws4.Range("B2").Value = "Estudo de como será feita a integração/consulta do banco de dados com o site através do Python3"
$ws4.Range("C2").Value = "x"
$ws4.Range("D2").Value = "x"
$ws4.Range("E2").Value = "Alta"

$ws4.Range("A3").Value = "Pesquisa sobre outras cidades"
$ws4.Range("B3").Value = "x"
$ws4.Range("C3").Value = "x"
$ws4.Range("D3").Value = "x"
$ws4.Range("E3").Value = "?"

# Row2 used to be styled s=20 on D2; the new layout reuses the common
# "x"-cell look (s=17/18) plus the Alta-priority look on E2/E3.
$ws3.Range("C3").Copy()
$ws4.Range("C2").PasteSpecial(-4122)
$ws4.Range("D2").PasteSpecial(-4122)
$ws4.Range("C3").PasteSpecial(-4122)
$ws4.Range("D3").PasteSpecial(-4122)
$ws4.Range("B3").PasteSpecial(-4122)
$ws3.Application.CutCopyMode = $false

$altaCell.Copy()
$ws4.Range("E2").PasteSpecial(-4122)
$ws4.Range("E3").PasteSpecial(-4122)
$ws3.Application.CutCopyMode = $false

$ws4.Rows.Item(1).RowHeight = 57
$ws4.Rows.Item(2).RowHeight = 84
$ws4.Rows.Item(3).RowHeight = 21

$ws4.PageSetup.PaperSize = 9
$ws4.PageSetup.Orientation = 1

$ws4.Range("E11").Select()

Write-Output "edit complete"
